$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Object" -> "Full ID"
$ws.Range("A1").Value = "Full ID"

# Rows 2-43: rename "*.bed" entries to "*.rds" (row 22 also fixes hg19->hg38 in filename)
$renames = @{
  2  = "ce10.Boyle.ce10-blacklist.v2.rds"
  3  = "ce10.Kundaje.ce10-blacklist.rds"
  4  = "ce11.Boyle.ce11-blacklist.v2.rds"
  5  = "danRer10.Domingues.blacklisted.rds"
  6  = "danRer10.Yang.Supplemental_Table_19.ChIP-seq_black_list_in_the_zebrafish_genome.rds"
  7  = "dm3.Boyle.dm3-blacklist.v2.rds"
  8  = "dm3.Kundaje.dm3-blacklist.rds"
  9  = "dm6.Boyle.dm6-blacklist.v2.rds"
  10 = "hg19.Bernstein.Mint_Blacklist_hg19.rds"
  11 = "hg19.Birney.wgEncodeDacMapabilityConsensusExcludable.rds"
  12 = "hg19.Boyle.hg19-blacklist.v2.rds"
  13 = "hg19.Crawford.wgEncodeDukeMapabilityRegionsExcludable.rds"
  14 = "hg19.Lareau.hg19_peaks.narrowPeak.rds"
  15 = "hg19.Lareau.hg19.full.blacklist.rds"
  16 = "hg19.Wold.hg19mitoblack.rds"
  17 = "hg19.Yeo.eCLIP_blacklistregions.hg19.rds"
  18 = "hg38.Bernstein.Mint_Blacklist_GRCh38.rds"
  19 = "hg38.Boyle.hg38-blacklist.v2.rds"
  20 = "hg38.Kundaje.GRCh38_unified_blacklist.rds"
  21 = "hg38.Kundaje.GRCh38.blacklist.rds"
  22 = "hg38.Lareau.hg38_peaks.narrowPeak.rds"
  23 = "hg38.Lareau.hg38.full.blacklist.rds"
  24 = "hg38.Reddy.wgEncodeDacMapabilityConsensusExcludable.hg38.rds"
  25 = "hg38.Wimberley.peakPass60Perc_sorted.rds"
  26 = "hg38.Wold.hg38mitoblack.rds"
  27 = "hg38.Yeo.eCLIP_blacklistregions.hg38liftover.bed.fixed.rds"
  28 = "mm10.Boyle.mm10-blacklist.v2.rds"
  29 = "mm10.Hardison.blacklist.full.rds"
  30 = "mm10.Hardison.psublacklist.mm10.rds"
  31 = "mm10.Kundaje.anshul.blacklist.mm10.rds"
  32 = "mm10.Kundaje.mm10.blacklist.rds"
  33 = "mm10.Lareau.mm10_peaks.narrowPeak.rds"
  34 = "mm10.Lareau.mm10.full.blacklist.rds"
  35 = "mm10.Wold.mm10mitoblack.rds"
  36 = "mm9.Lareau.mm9_peaks.narrowPeak.rds"
  37 = "mm9.Lareau.mm9.full.blacklist.rds"
  38 = "mm9.Wold.mm9mitoblack.rds"
  39 = "T2T.excluderanges.excludable.rds"
  40 = "T2T.Lareau.chm13v2.0_peaks.narrowPeak.rds"
  41 = "TAIR10.Klasfeld.arabidopsis_blacklist_20inputs.rds"
  42 = "TAIR10.Klasfeld.arabidopsis_greenscreen_20inputs.rds"
  43 = "TAIR10.Wimberley.predicted_excluded_list_sorted_0.6.rds"
}

foreach ($r in $renames.Keys) {
  $ws.Cells.Item($r, 1).Value = $renames[$r]
}

# Row 51: correct mismatched hg38 centromere entry
$ws.Range("A51").Value = "hg38.UCSC.centromere.bed"
$ws.Range("B51").Value = "hg38"
$ws.Range("D51").Value = "11 : 3,391 : 30,000,001"

# New row 68: T2T UCSC hgUnique entry
$ws.Range("A68").Value = "T2T.UCSC.hgUnique.hg38.bed"
$ws.Range("B68").Value = "T2T"
$ws.Range("C68").Value = 615
$ws.Range("D68").Value = "2 : 15,829 : 29,694,330"
$ws.Range("E68").Value = "M"
$ws.Range("F68").Value = 2022
$ws.Range("G68").Value = "https://hgdownload.soe.ucsc.edu/hubs/GCA/009/914/755/GCA_009914755.4/bbi/GCA_009914755.4_T2T-CHM13v2.0.hgUnique/hgUnique.hg38.bb"

Write-Output "edit applied"
